$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "10/17/2025"
$ws.Range("A10").ClearFormats()
$ws.Range("B10").Value = 463.9639999999999
$ws.Range("C10").Value = 0.1077669819210111
$ws.Range("D10").Value = 25
